$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '86.620.06'
$ws.Range("E2").Value = '  +3.33%  '

$ws.Range("D3").Value = '3.279.10'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.43'
$ws.Range("E5").Value = '  -3.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.02'
$ws.Range("E6").Value = '  +0.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.368'
$ws.Range("E7").Value = '  +20.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.697'
$ws.Range("E8").Value = '  +18.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '3.275.30'
$ws.Range("E10").Value = '  +0.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  -5.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.179'
$ws.Range("E12").Value = '  +7.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  -7.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.10'
$ws.Range("E14").Value = '  +4.68%  '

$ws.Range("D15").Value = '3.875.77'
$ws.Range("E15").Value = '  +0.97%  '

$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").Value = '86.587.53'
$ws.Range("E17").Value = '  +3.96%  '

$ws.Range("D18").Value = '3.266.03'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.02'
$ws.Range("E19").Value = '  -2.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.04'
$ws.Range("E20").Value = '  -5.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.84'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.90'
$ws.Range("E22").Value = '  -1.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.32'
$ws.Range("E23").Value = '  +1.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.27'
$ws.Range("E24").Value = '  -2.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.50'
$ws.Range("E25").Value = '  +5.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.11'
$ws.Range("E26").Value = '  -1.79%  '

$ws.Range("D27").Value = '3.444.94'
$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("E28").Value = '  -3.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000129'
$ws.Range("E29").Value = '  +4.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("E31").Value = '  +13.94%  '

$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.81'
$ws.Range("E33").Value = '  -3.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '543.10'
$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("E35").Value = '  -3.90%  '

$ws.Range("E36").Value = '  -2.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.94'
$ws.Range("E37").Value = '  +12.24%  '

$ws.Range("E38").Value = '  -10.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.49'
$ws.Range("E39").Value = '  -2.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.58'
$ws.Range("E41").Value = '  +3.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.394'
$ws.Range("E42").Value = '  -3.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").Value = '  -1.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.92'
$ws.Range("E44").Value = '  -3.18%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '157.43'
$ws.Range("E46").Value = '  -1.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '179.41'
$ws.Range("E47").Value = '  -4.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.43'
$ws.Range("E48").Value = '  -1.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").Value = '  -0.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.24'
$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.123'
$ws.Range("E51").Value = '  +11.80%  '
